$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update report title (month changed from October 2016 to November 2016)
$ws.Range("A2").Value = "Electric Utilities by Census Division and State, November 2016 (Continued)"

# Update Relative Standard Error data cells
$ws.Range("E4").Value = 58
$ws.Range("I4").Value = 14
$ws.Range("E5").Value = 138
$ws.Range("F5").Value = 138
$ws.Range("I5").Value = 103
$ws.Range("I6").Value = 405
$ws.Range("E7").Value = 64
$ws.Range("F7").Value = 34
$ws.Range("I7").Value = 57
$ws.Range("I8").Value = 11
$ws.Range("I9").Value = 32
$ws.Range("I10").Value = 24
$ws.Range("E11").Value = 28
$ws.Range("F11").Value = 28
$ws.Range("I11").Value = 3
$ws.Range("E12").Value = 28
$ws.Range("F12").Value = 28
$ws.Range("I12").Value = 37
$ws.Range("I13").Value = 3
$ws.Range("I14").Value = 156
$ws.Range("E15").Value = 27
$ws.Range("F16").Value = 55
$ws.Range("I16").Value = 2
$ws.Range("E17").Value = 34
$ws.Range("F17").Value = 18
$ws.Range("E18").Value = 58
$ws.Range("F18").Value = 1
$ws.Range("E19").Value = 79
$ws.Range("F19").Value = 54
$ws.Range("F21").Value = 0.39
$ws.Range("H21").Value = 9
$ws.Range("F22").Value = 0.25
$ws.Range("I23").Value = 1
$ws.Range("F25").Value = 50
$ws.Range("F26").Value = 9
$ws.Range("H27").Value = 54
$ws.Range("F28").Value = 0.39
$ws.Range("I28").Value = 2
$ws.Range("E29").Value = 8
$ws.Range("I29").Value = 0.17
$ws.Range("E30").Value = 97
$ws.Range("F30").Value = 97
$ws.Range("I30").Value = 124
$ws.Range("E31").Value = 11
$ws.Range("F31").Value = 10
$ws.Range("I31").Value = 0.37
$ws.Range("E32").Value = 4
$ws.Range("F32").Value = 4
$ws.Range("I32").Value = 0.28999999999999998
$ws.Range("E33").Value = 84
$ws.Range("F33").Value = 84
$ws.Range("I33").Value = 44
$ws.Range("E34").Value = 22
$ws.Range("F34").Value = 22
$ws.Range("I34").Value = 0.39
$ws.Range("F35").Value = 8
$ws.Range("I35").Value = 0.39
$ws.Range("E36").Value = 118
$ws.Range("F36").Value = 2
$ws.Range("I36").Value = 0.31
$ws.Range("F38").Value = 33
$ws.Range("I38").Value = 0.46
$ws.Range("F40").Value = 33
$ws.Range("I43").Value = 0.42
$ws.Range("I45").Value = 0.27
$ws.Range("E48").Value = 8
$ws.Range("H48").Value = 2257
$ws.Range("I48").Value = 1
$ws.Range("E49").Value = 9
$ws.Range("F49").Value = 9
$ws.Range("I49").Value = 0.19
$ws.Range("F50").Value = 5
$ws.Range("F51").Value = 100
$ws.Range("I51").Value = 10
$ws.Range("I52").Value = 5
$ws.Range("E53").Value = 53
$ws.Range("F53").Value = 53
$ws.Range("I53").Value = 0.05
$ws.Range("E54").Value = 18
$ws.Range("F54").Value = 18
$ws.Range("H54").Value = 2158
$ws.Range("E57").Value = 12
$ws.Range("E58").Value = 12
$ws.Range("F58").Value = 3
$ws.Range("C59").Value = 222
$ws.Range("E59").Value = 98
$ws.Range("F59").Value = 8
$ws.Range("I59").Value = 2
$ws.Range("E61").Value = 41
$ws.Range("F61").Value = 18
$ws.Range("I61").Value = 5
$ws.Range("F62").Value = 29
$ws.Range("E63").Value = 41
$ws.Range("F63").Value = 17
$ws.Range("I63").Value = 6
$ws.Range("F64").Value = 0.41
$ws.Range("H64").Value = 7
$ws.Range("I64").Value = 0.19
